$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.406.40'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.693.20'
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = "'219.22"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = "'0.5486"
$ws.Range("E6").Value = '  +4.53%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = "'0.2732"
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").Value = "'0.06470"
$ws.Range("D10").Value = "'22.03"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = "'0.07679"
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("D12").Value = '1.690.14'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = "'4.554"
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").Value = "'0.5849"
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = "'0.000008388"
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = "'65.45"
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '26.435.23'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = "'4.953"
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = "'10.98"
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = "'192.53"
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").Value = "'6.260"
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = "'149.46"
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("E25").Value = '  +7.74%  '
$ws.Range("D26").Value = "'7.913"
$ws.Range("E26").Value = '  +3.19%  '
$ws.Range("D27").Value = "'15.77"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("E28").Value = '  -5.38%  '
$ws.Range("D29").Value = "'1.396"
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = "'1.332"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = "'3.603"
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("D35").Value = "'0.6154"
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("D36").Value = "'2.409"
$ws.Range("E36").Value = '  +0.75%  '
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = "'6.209"
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("D39").Value = '1.120.50'
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("D40").Value = "'0.01639"
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").Value = "'0.8832"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = "'101.97"
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D44").Value = '1.843.53'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'57.55"
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").Value = "'8.219"
$ws.Range("E47").Value = '  +0.65%  '
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = "'0.05284"
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.4305"
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = "'6.103"
$ws.Range("E51").Value = '  +0.71%  '
